$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (single-decimal numeric-looking strings) are explicitly formatted as Text
# first, matching how a user would force text entry in real Excel.

$ws.Range('D2').Value = '64.324.85'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.487.50'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.87'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('D7').Value = '3.488.50'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('E12').Value = '  -2.39%  '
$ws.Range('D13').Value = '4.081.15'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = '3.487.50'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').Value = '64.352.42'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.14'
$ws.Range('E18').Value = '  -9.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.96'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.67'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  -3.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '385.70'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.567'
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('D24').Value = '3.624.73'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.01'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  +4.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.56'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.45'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.21'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').Value = '3.507.23'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +1.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '23.40'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.30'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.84'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '162.38'
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0779'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.802'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.72'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.73'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('D50').Value = '2.469.90'
$ws.Range('E50').Value = '  +2.02%  '
$ws.Range('E51').Value = '  -1.65%  '
